$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 237.1239
$ws.Range("C3").Value = 241.4432
$ws.Range("C4").Value = 241.5905
$ws.Range("C5").Value = 238.5242
$ws.Range("C6").Value = 240.5607
$ws.Range("C7").Value = 241.3397
$ws.Range("C8").Value = 241.9323
$ws.Range("C9").Value = 244.4068
$ws.Range("C10").Value = 241.0366
$ws.Range("C11").Value = 232.9335
$ws.Range("C12").Value = 227.9097
$ws.Range("C13").Value = 225.2887
$ws.Range("C14").Value = 230.8064
$ws.Range("C15").Value = 235.107
$ws.Range("C16").Value = 232.6998
$ws.Range("C17").Value = 223.5023
$ws.Range("C18").Value = 228.8661
$ws.Range("C19").Value = 239.0934
$ws.Range("C20").Value = 238.2435
$ws.Range("C21").Value = 234.8609
$ws.Range("C22").Value = 225.5632
$ws.Range("C23").Value = 221.695
$ws.Range("C24").Value = 223.4547
